$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("P2").Value = 2.12
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 2.1
$ws.Range("W2").Value = 2.38

# Row 4 updates
$ws.Range("F4").Value = 8.6
$ws.Range("H4").Value = 1.42
$ws.Range("I4").Value = 1.44
$ws.Range("K4").Value = 5.4
$ws.Range("Q4").Value = 1.88
$ws.Range("U4").Value = 1.77
